$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 5.198314948599565
$ws.Cells.Item(2, 4).Value = 7.897489847046764
$ws.Cells.Item(2, 5).Value = 8.93155618175153
$ws.Cells.Item(2, 6).Value = 70.52646320639889
$ws.Cells.Item(2, 7).Value = 3.906213214409895
$ws.Cells.Item(2, 9).Value = 53.84651163312628
$ws.Cells.Item(2, 10).Value = 9.602474674615609
$ws.Cells.Item(2, 12).Value = 8.540587365727752
$ws.Cells.Item(2, 13).Value = 66.37944984824257
$ws.Cells.Item(3, 3).Value = 5.052345454364529
$ws.Cells.Item(3, 4).Value = 7.863146849621405
$ws.Cells.Item(3, 5).Value = 8.657496921354131
$ws.Cells.Item(3, 6).Value = 71.45863115067127
$ws.Cells.Item(3, 7).Value = 3.916898525676623
$ws.Cells.Item(3, 9).Value = 54.54113175324464
$ws.Cells.Item(3, 10).Value = 9.659684814723612
$ws.Cells.Item(3, 12).Value = 8.501297239242664
$ws.Cells.Item(3, 13).Value = 63.71465339034275
$ws.Cells.Item(4, 3).Value = 4.963110341199064
$ws.Cells.Item(4, 4).Value = 7.846607481981634
$ws.Cells.Item(4, 5).Value = 8.485422718989593
$ws.Cells.Item(4, 6).Value = 72.06856612094353
$ws.Cells.Item(4, 7).Value = 3.923721037419392
$ws.Cells.Item(4, 9).Value = 54.9948945358271
$ws.Cells.Item(4, 10).Value = 9.697523365489671
$ws.Cells.Item(4, 12).Value = 8.478792269695928
$ws.Cells.Item(4, 13).Value = 62.02019375725468
$ws.Cells.Item(5, 3).Value = 4.92689865809255
$ws.Cells.Item(5, 4).Value = 7.840996528909931
$ws.Cells.Item(5, 5).Value = 8.414425498958318
$ws.Cells.Item(5, 6).Value = 72.32642482199192
$ws.Cells.Item(5, 7).Value = 3.926567825991617
$ws.Cells.Item(5, 9).Value = 55.18656710032047
$ws.Cells.Item(5, 10).Value = 9.713619044758612
$ws.Cells.Item(5, 12).Value = 8.470029273197541
$ws.Cells.Item(5, 13).Value = 61.31571528457734
$ws.Cells.Item(6, 3).Value = 4.920896552811271
$ws.Cells.Item(6, 4).Value = 7.840132558985211
$ws.Cells.Item(6, 5).Value = 8.402586060131952
$ws.Cells.Item(6, 6).Value = 72.36980004198553
$ws.Cells.Item(6, 7).Value = 3.927044574819644
$ws.Cells.Item(6, 9).Value = 55.2187998544518
$ws.Cells.Item(6, 10).Value = 9.716332402735102
$ws.Cells.Item(6, 12).Value = 8.468598843414913
$ws.Cells.Item(6, 13).Value = 61.19791498867772
$ws.Cells.Item(7, 3).Value = 4.962621288099344
$ws.Cells.Item(7, 4).Value = 7.846527260136887
$ws.Cells.Item(7, 5).Value = 8.484468659849419
$ws.Cells.Item(7, 6).Value = 72.07200620037626
$ws.Cells.Item(7, 7).Value = 3.923759159720765
$ws.Cells.Item(7, 9).Value = 54.99745224629395
$ws.Cells.Item(7, 10).Value = 9.697737706773076
$ws.Cells.Item(7, 12).Value = 8.478672436135405
$ws.Cells.Item(7, 13).Value = 62.01074852110496
$ws.Cells.Item(8, 3).Value = 5.147935400058694
$ws.Cells.Item(8, 4).Value = 7.884696463454653
$ws.Cells.Item(8, 5).Value = 8.837895274092347
$ws.Cells.Item(8, 6).Value = 70.83998488047506
$ws.Cells.Item(8, 7).Value = 3.90984365954376
$ws.Cells.Item(8, 9).Value = 54.08029893371777
$ws.Cells.Item(8, 10).Value = 9.621634491465652
$ws.Cells.Item(8, 12).Value = 8.52670453121611
$ws.Cells.Item(8, 13).Value = 65.47309843164243
$ws.Cells.Item(9, 3).Value = 5.51217719704173
$ws.Cells.Item(9, 4).Value = 7.996176813516407
$ws.Cells.Item(9, 5).Value = 9.497651815102051
$ws.Cells.Item(9, 6).Value = 68.72919080384573
$ws.Cells.Item(9, 7).Value = 3.884596223187341
$ws.Cells.Item(9, 9).Value = 52.50269999817293
$ws.Cells.Item(9, 10).Value = 9.494183864829669
$ws.Cells.Item(9, 12).Value = 8.633710855245027
$ws.Cells.Item(9, 13).Value = 71.775380337377
$ws.Cells.Item(10, 3).Value = 5.777485499830691
$ws.Cells.Item(10, 4).Value = 8.101023762939075
$ws.Cells.Item(10, 5).Value = 9.958354598976525
$ws.Cells.Item(10, 6).Value = 67.37498966563457
$ws.Cells.Item(10, 7).Value = 3.867239571461962
$ws.Cells.Item(10, 9).Value = 51.4852720846751
$ws.Cells.Item(10, 10).Value = 9.414242840509178
$ws.Cells.Item(10, 12).Value = 8.720101443480219
$ws.Cells.Item(10, 13).Value = 76.08102422280987
$ws.Cells.Item(11, 3).Value = 5.89715619657213
$ws.Cells.Item(11, 4).Value = 8.153792122913559
$ws.Cells.Item(11, 5).Value = 10.16208039864228
$ws.Cells.Item(11, 6).Value = 66.80420437405306
$ws.Cells.Item(11, 7).Value = 3.859590807414081
$ws.Cells.Item(11, 9).Value = 51.0548811490119
$ws.Cells.Item(11, 10).Value = 9.380951529969979
$ws.Cells.Item(11, 12).Value = 8.761081897178613
$ws.Cells.Item(11, 13).Value = 77.96499559874977
$ws.Cells.Item(12, 3).Value = 5.942287124482649
$ws.Cells.Item(12, 4).Value = 8.174507903189186
$ws.Cells.Item(12, 5).Value = 10.23833982501721
$ws.Cells.Item(12, 6).Value = 66.59480907937551
$ws.Cells.Item(12, 7).Value = 3.85672892668174
$ws.Cells.Item(12, 9).Value = 50.89672743809314
$ws.Cells.Item(12, 10).Value = 9.368796462227612
$ws.Cells.Item(12, 12).Value = 8.7768408801494
$ws.Cells.Item(12, 13).Value = 78.66739442449371
$ws.Cells.Item(13, 3).Value = 5.932576214698847
$ws.Cells.Item(13, 4).Value = 8.170013718318614
$ws.Cells.Item(13, 5).Value = 10.22195611129865
$ws.Cells.Item(13, 6).Value = 66.63960188133322
$ws.Cells.Item(13, 7).Value = 3.85734376342523
$ws.Cells.Item(13, 9).Value = 50.9305713242622
$ws.Cells.Item(13, 10).Value = 9.371394026998479
$ws.Cells.Item(13, 12).Value = 8.773436218898354
$ws.Cells.Item(13, 13).Value = 78.51661424811897
$ws.Cells.Item(14, 3).Value = 5.900873078331341
$ws.Cells.Item(14, 4).Value = 8.15548170416252
$ws.Cells.Item(14, 5).Value = 10.16837232699514
$ws.Cells.Item(14, 6).Value = 66.78684058438226
$ws.Cells.Item(14, 7).Value = 3.859354672163567
$ws.Cells.Item(14, 9).Value = 51.04177205776024
$ws.Cells.Item(14, 10).Value = 9.379942416488667
$ws.Cells.Item(14, 12).Value = 8.762373582977196
$ws.Cells.Item(14, 13).Value = 78.02300463488999
$ws.Cells.Item(15, 3).Value = 5.881428705350233
$ws.Cells.Item(15, 4).Value = 8.146676046138026
$ws.Cells.Item(15, 5).Value = 10.13543392984864
$ws.Cells.Item(15, 6).Value = 66.87791495156651
$ws.Cells.Item(15, 7).Value = 3.860590880834189
$ws.Cells.Item(15, 9).Value = 51.11051921220986
$ws.Cells.Item(15, 10).Value = 9.385237663794896
$ws.Cells.Item(15, 12).Value = 8.755628696387985
$ws.Cells.Item(15, 13).Value = 77.71921200541398
$ws.Cells.Item(16, 3).Value = 5.769640784662915
$ws.Cells.Item(16, 4).Value = 8.09767769509836
$ws.Cells.Item(16, 5).Value = 9.944918853106657
$ws.Cells.Item(16, 6).Value = 67.41322257252543
$ws.Cells.Item(16, 7).Value = 3.867744322709319
$ws.Cells.Item(16, 9).Value = 51.51406540582175
$ws.Cells.Item(16, 10).Value = 9.416481167571055
$ws.Cells.Item(16, 12).Value = 8.717456999987983
$ws.Cells.Item(16, 13).Value = 75.95637405744843
$ws.Cells.Item(17, 3).Value = 5.700772277676957
$ws.Cells.Item(17, 4).Value = 8.068921558010032
$ws.Cells.Item(17, 5).Value = 9.826510998177765
$ws.Cells.Item(17, 6).Value = 67.75337153489302
$ws.Cells.Item(17, 7).Value = 3.87219528617879
$ws.Cells.Item(17, 9).Value = 51.77004735155219
$ws.Cells.Item(17, 10).Value = 9.43644199984379
$ws.Cells.Item(17, 12).Value = 8.694469596504877
$ws.Cells.Item(17, 13).Value = 74.85558251359457
$ws.Cells.Item(18, 3).Value = 5.661067076323194
$ws.Cells.Item(18, 4).Value = 8.052859154946956
$ws.Cells.Item(18, 5).Value = 9.757857610732334
$ws.Cells.Item(18, 6).Value = 67.95325947530056
$ws.Cells.Item(18, 7).Value = 3.874778658395972
$ws.Cells.Item(18, 9).Value = 51.92032457987344
$ws.Cells.Item(18, 10).Value = 9.448211690425625
$ws.Cells.Item(18, 12).Value = 8.681406134259197
$ws.Cells.Item(18, 13).Value = 74.21541917031225
$ws.Cells.Item(19, 3).Value = 5.647608716875911
$ws.Cells.Item(19, 4).Value = 8.047502536402273
$ws.Cells.Item(19, 5).Value = 9.734520079940268
$ws.Cells.Item(19, 6).Value = 68.02165982376279
$ws.Cells.Item(19, 7).Value = 3.875657373448045
$ws.Cells.Item(19, 9).Value = 51.97172385481477
$ws.Cells.Item(19, 10).Value = 9.452246019207342
$ws.Cells.Item(19, 12).Value = 8.677010291792323
$ws.Cells.Item(19, 13).Value = 73.99747516509959
$ws.Cells.Item(20, 3).Value = 5.708113458335221
$ws.Cells.Item(20, 4).Value = 8.071933257253223
$ws.Cells.Item(20, 5).Value = 9.839172817845222
$ws.Cells.Item(20, 6).Value = 67.71672124563376
$ws.Cells.Item(20, 7).Value = 3.871719069310795
$ws.Cells.Item(20, 9).Value = 51.7424815563537
$ws.Cells.Item(20, 10).Value = 9.43428718405881
$ws.Cells.Item(20, 12).Value = 8.696900271543841
$ws.Cells.Item(20, 13).Value = 74.97349211791939
$ws.Cells.Item(21, 3).Value = 5.910190391930864
$ws.Cells.Item(21, 4).Value = 8.159730176977192
$ws.Cells.Item(21, 5).Value = 10.18413559370571
$ws.Cells.Item(21, 6).Value = 66.74340782057465
$ws.Cells.Item(21, 7).Value = 3.85876308995543
$ws.Cells.Item(21, 9).Value = 51.00897738073229
$ws.Cells.Item(21, 10).Value = 9.377419210294079
$ws.Cells.Item(21, 12).Value = 8.765616426892485
$ws.Cells.Item(21, 13).Value = 78.16829077682893
$ws.Cells.Item(22, 3).Value = 6.041161234934413
$ws.Cells.Item(22, 4).Value = 8.221386906612626
$ws.Cells.Item(22, 5).Value = 10.40440102202188
$ws.Cells.Item(22, 6).Value = 66.14675483848316
$ws.Cells.Item(22, 7).Value = 3.850496460311455
$ws.Cells.Item(22, 9).Value = 50.55780420461313
$ws.Cells.Item(22, 10).Value = 9.34289012463687
$ws.Cells.Item(22, 12).Value = 8.811928198812209
$ws.Cells.Item(22, 13).Value = 80.19195154883907
$ws.Cells.Item(23, 3).Value = 5.971372009034657
$ws.Cells.Item(23, 4).Value = 8.188087361768325
$ws.Cells.Item(23, 5).Value = 10.28732943672109
$ws.Cells.Item(23, 6).Value = 66.46150354342328
$ws.Cells.Item(23, 7).Value = 3.854890465791032
$ws.Cells.Item(23, 9).Value = 50.79596546366209
$ws.Cells.Item(23, 10).Value = 9.361074267261479
$ws.Cells.Item(23, 12).Value = 8.787082824329048
$ws.Cells.Item(23, 13).Value = 79.11784934705555
$ws.Cells.Item(24, 3).Value = 5.704794855708111
$ws.Cells.Item(24, 4).Value = 8.070570204771069
$ws.Cells.Item(24, 5).Value = 9.833450209498974
$ws.Cells.Item(24, 6).Value = 67.73327736698623
$ws.Cells.Item(24, 7).Value = 3.871934290715874
$ws.Cells.Item(24, 9).Value = 51.75493437804053
$ws.Cells.Item(24, 10).Value = 9.435260461424368
$ws.Cells.Item(24, 12).Value = 8.695800888882426
$ws.Cells.Item(24, 13).Value = 74.9202079130644
$ws.Cells.Item(25, 3).Value = 5.413853175747871
$ws.Cells.Item(25, 4).Value = 7.962025755682467
$ws.Cells.Item(25, 5).Value = 9.323147040927061
$ws.Cells.Item(25, 6).Value = 69.26651003094373
$ws.Cells.Item(25, 7).Value = 3.891213044456442
$ws.Cells.Item(25, 9).Value = 52.90514984038943
$ws.Cells.Item(25, 10).Value = 9.526289519698926
$ws.Cells.Item(25, 12).Value = 8.603394601576939
$ws.Cells.Item(25, 13).Value = 70.12576973934311
